# Apply the commit's edits to the workbook:
#  - Sheet2!B3 (Pax 1 Firstname) "Tamil" -> "Tamil Yadav"
#  - Sheet2!B12 (Pax 2 Firstname) "Tamil" -> "Lingeswar"
#  - Selection/active sheet moves from Sheet1 (at B5) to Sheet2 (at B3),
#    leaving Sheet1's lingering selection at B2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update the pax names on Sheet2
$ws2.Range("B3").Value = "Tamil Yadav"
$ws2.Range("B12").Value = "Lingeswar"

# Leave Sheet1's selection parked at B2 (no longer the active tab)
$ws1.Activate()
$ws1.Range("B2").Select()

# Make Sheet2 the active tab with B3 selected
$ws2.Activate()
$ws2.Range("B3").Select()
